$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.205562353134155
$ws.Range("B1").Value = 3.28863787651062
$ws.Range("C1").Value = 2.774307727813721
$ws.Range("D1").Value = 2.102351665496826
$ws.Range("E1").Value = 1.248859643936157
